$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,("AG8700PRO3", "22WNBPEYG6K1BU35", 5.73, 6, "trial")
    ,("AS1633PRO3", "22WNBPEYG6K1BU35", 4.86, 6, "trial")
    ,("AS1820PRO3", "22WNBPEYG6K1BU35", 3.99, 6, "trial")
    ,("AS1868PRO3", "22WNBPEYG6K1BU35", 4.34, 6, "trial")
    ,("B2401PWU", "22WNBPEYG6K1BU35", 4.51, 6, "trial")
    ,("B2620PWU", "22WNBPEYG6K1BU35", 5.21, 6, "trial")
    ,("DKB255PRO3", "22WNBPEYG6K1BU35", 3.99, 6, "trial")
    ,("DKB335PRO3", "22WNBPEYG6K1BU35", 4.16, 6, "trial")
    ,("DKB360PRO3", "22WNBPEYG6K1BU35", 4.51, 6, "trial")
    ,("K9606VIP3", "22WNBPEYG6K1BU35", 4.69, 6, "trial")
    ,("MG580PWU", "22WNBPEYG6K1BU35", 4.69, 6, "trial")
    ,("MG593PWU", "22WNBPEYG6K1BU35", 3.99, 6, "trial")
    ,("MG652PWU", "22WNBPEYG6K1BU35", 4.16, 6, "trial")
    ,("NK506VIP3", "22WNBPEYG6K1BU35", 4.34, 6, "trial")
    ,("NK508", "22WNBPEYG6K1BU35", 3.81, 6, "trial")
    ,("NK511VIP3", "22WNBPEYG6K1BU35", 5.38, 6, "trial")
    ,("NK520VIP3", "22WNBPEYG6K1BU35", 4.69, 6, "trial")
    ,("NK555VIP3", "22WNBPEYG6K1BU35", 5.64, 6, "trial")
    ,("NS75VIP3", "22WNBPEYG6K1BU35", 4.86, 6, "trial")
    ,("NS80VIP3", "22WNBPEYG6K1BU35", 5.38, 6, "trial")
    ,("NS88VIP3", "22WNBPEYG6K1BU35", 4.51, 6, "trial")
    ,("NS90PRO2", "22WNBPEYG6K1BU35", 3.99, 6, "trial")
    ,("NS91VIP2", "22WNBPEYG6K1BU35", 4.86, 6, "trial")
    ,("NS91VIP3", "22WNBPEYG6K1BU35", 4.69, 6, "trial")
    ,("NS95VIP2", "22WNBPEYG6K1BU35", 5.21, 6, "trial")
    ,("P3707VYH", "22WNBPEYG6K1BU35", 4.86, 6, "trial")
    ,("P3845VYHR", "22WNBPEYG6K1BU35", 5.38, 6, "trial")
    ,("P3858PWU", "22WNBPEYG6K1BU35", 5.21, 6, "trial")
    ,("P3889RR", "22WNBPEYG6K1BU35", 4.86, 6, "trial")
    ,("P3898", "22WNBPEYG6K1BU35", 5.21, 6, "trial")
    ,("STATUSVIP3", "22WNBPEYG6K1BU35", 3.99, 6, "trial")
    ,("SXC2251TG.0", "22WNBPEYG6K1BU35", 4.51, 6, "trial")
    ,("SXC2251ZL.0", "22WNBPEYG6K1BU35", 4.34, 6, "trial")
    ,("SXC2320VIP3", "22WNBPEYG6K1BU35", 4.86, 6, "trial")
    ,("SXD2621ZL.0", "22WNBPEYG6K1BU35", 3.99, 6, "trial")
    ,("SYC7131ZL.0", "22WNBPEYG6K1BU35", 3.99, 6, "trial")
    ,("SYC7441ZL.0", "22WNBPEYG6K1BU35", 5.04, 6, "trial")
    ,("SZC4141ZL.0", "22WNBPEYG6K1BU35", 4.86, 6, "trial")
    ,("SZD5031ZL.0", "22WNBPEYG6K1BU35", 3.99, 6, "trial")
    ,("SZD5211ZL.0", "22WNBPEYG6K1BU35", 4.86, 6, "trial")
    ,("SZE6111TG.0", "22WNBPEYG6K1BU35", 3.99, 6, "trial")
    ,("SZE6111VIP3", "22WNBPEYG6K1BU35", 3.99, 6, "trial")
    ,("SZF6101VIP3", "22WNBPEYG6K1BU35", 3.64, 6, "trial")
    ,("SZF6101WT.0", "22WNBPEYG6K1BU35", 3.64, 6, "trial")
)

$startRow = 1947
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

